# Starting exploratory analysis for the dataset:
# add a new summary row (row 6) that expresses each indicator's "share"
# (row 5) as a percentage of its total (row 2), widen the two columns
# that needed extra room, and move the selection onto the freshly
# populated cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: for every data column (B..AT) compute <share>/<total> * 100,
# based on the existing row 2 (totals) and row 5 (share counts).
$ws.Range("B6:AT6").Formula = "=B5/B2 * 100"

# Columns D and E need to be a bit wider to comfortably show the new
# percentage values instead of the bestFit widths used for the old
# content.
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666

# Reflect the work-in-progress selection on the new row.
[void]$ws.Range("E4:E6").Select()
